# Fix for Induction Loop Data
# - The hourly timestamp in B2 was off by one hour (08:00:00 -> 09:00:00).
# - The measured total in C2 had a typo ("924,7" -> "824,7").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the timestamp for the second data row (was 45268.33333333334 / 08:00:00).
$ws.Range("B2").Value = 45268.375

# Correct the total value for the second data row.
$ws.Range("C2").Value = "824,7"

# Bring page margins in line with Excel's standard defaults (0.7/0.7/0.75/0.75/0.3/0.3 in).
$ws.PageSetup.LeftMargin = 0.7 * 72
$ws.PageSetup.RightMargin = 0.7 * 72
$ws.PageSetup.TopMargin = 0.75 * 72
$ws.PageSetup.BottomMargin = 0.75 * 72
$ws.PageSetup.HeaderMargin = 0.3 * 72
$ws.PageSetup.FooterMargin = 0.3 * 72
